# Edit: split "Xem tat ca sach ... UCCL-01" into two runs:
#   run 1: "Xem tất cả sách"
#   run 2: " – UCN-05"   (en dash, leading space, xml:space="preserve")
#
# The original text lives in a single run/paragraph:
#   "Xem tất cả sách                                      UCCL-01"
# We first trim it down (via wildcard Find & Replace) to just the title,
# then insert the new suffix as a genuinely separate run (forcing a run
# boundary by toggling Bold on/off around the InsertAfter call, which the
# engine preserves as two runs with identical resulting formatting).

$d = $word.ActiveDocument

# Step 1: collapse "Xem tất cả sách<spaces>UCCL-01" down to just the title.
$find1 = $d.Content
$found1 = $find1.Find.Execute(
    "Xem tất cả sách*UCCL-01",  # FindText (wildcard)
    $false,                      # MatchCase
    $false,                      # MatchWholeWord
    $true,                       # MatchWildcards
    $false,                      # MatchSoundsLike
    $false,                      # MatchAllWordForms
    $true,                       # Forward
    1,                           # Wrap (wdFindContinue)
    $false,                      # Format
    "Xem tất cả sách",           # ReplaceWith
    2                            # Replace (wdReplaceAll)
)
if (-not $found1) {
    throw "Could not find the 'Xem tat ca sach ... UCCL-01' text to update."
}

# Step 2: locate the (now shortened) title text again to get its end position.
$find2 = $d.Content
$found2 = $find2.Find.Execute(
    "Xem tất cả sách",
    $false, $false, $false, $false, $false,
    $true, 1, $false, "", 0
)
if (-not $found2) {
    throw "Could not re-locate 'Xem tat ca sach' after trimming."
}

$insertPos = $find2.End
$insertRange = $d.Range($insertPos, $insertPos)

# Insert the new suffix as its own run. Toggling Bold on then back off
# around the insertion forces a distinct run boundary instead of merging
# back into the preceding run (even though the final formatting matches).
$insertRange.InsertAfter(" – UCN-05")
$suffixRange = $d.Range($insertPos, $insertPos + 9)
$suffixRange.Bold = 1
$suffixRange.Bold = 0
